$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at D; this shifts the existing D:Q columns
# (in every row) one position to the right, into E:R.
$ws.Columns("D").Insert()

# New header label for the freshly inserted column D.
$ws.Cells.Item(1, 4).Value = "Unnamed: 0.3"

# The data rows (2-20) already carry the duplicated "Unnamed: 0.x" index
# value in (the now-shifted) column E; copy that same value into the new
# column D so it matches its neighbours, exactly like the other repeated
# index columns. Row 21 is the trailing blank row, so its new column D
# cell just mirrors the (already blank) neighbouring cells.
for ($r = 2; $r -le 21; $r++) {
    $srcCell = $ws.Cells.Item($r, 5)
    $val = $srcCell.Value()
    $ws.Cells.Item($r, 4).Value = $val
}
